$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The "Team Members" table on slide 1 (graphicFrame shape 5).
$sh = $s.Shapes.Item(5)
$tbl = $sh.Table

# Row 2 / Column 1 held four team-member lines; trim it down to just
# "S SUJAN", dropping "P KIRAN", "S MUNI AJEY" and the trailing blank line.
$cell = $tbl.Cell(2, 1)
$cell.Shape.TextFrame.TextRange.Text = "S SUJAN             :RA2311003020537"

# PowerPoint re-flowed the table once the text shrank -- match the
# resulting row heights (EMU -> points, 12700 EMU per point).
$tbl.Rows.Item(2).Height = 1216153 / 12700.0
$tbl.Rows.Item(1).Height = 319288 / 12700.0
